# Update Harmonizing SXI-DF.xlsx: add "Prozess" sheet and two rows to "Zielstellung"
$wb = $excel.ActiveWorkbook

# --- 1. Add a new worksheet, rename it to "Prozess", and move it after "Datenstruktur" ---
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "Prozess"

$wsDaten = $wb.Worksheets.Item("Datenstruktur")
$wsNew.Move($null, $wsDaten)

# Re-fetch the worksheet object by name now that the sheet order has changed,
# since sheet references can become stale (positional) after a Move().
$wsProzess = $wb.Worksheets.Item("Prozess")

$wsProzess.Range("B4").Value = "Download SXI Logik als XML Dateien möglich"
$wsProzess.Range("B6").Value = "NVARCHAR muss in SXI eingestellt werden"
$wsProzess.Range("C5").Value = "Upload von geänderten Dateien pro Kunde mit Versionierung"
$wsProzess.Range("C7").Value = "Besser create Script für Datenbank für flexible Tabellendefinition ?"
$wsProzess.Range("B9").Value = "Webserver Yeti läuft als Dienst auf Port 8089"
$wsProzess.Range("B10").Value = "Vadin generiert HTML Seiten"
$wsProzess.Range("B11").Value = "Jenkins, Nexus, Maeven für Build -> auf Azure verfügbar - in VM, Container, Service ?"
$wsProzess.Range("B13").Value = "Repo im Stash auf Azure"

# --- 2. Add two rows to the "Zielstellung" sheet ---
$wsZiel = $wb.Worksheets.Item("Zielstellung")
$wsZiel.Range("B50").Value = "Beispiel Greif Velux"
$wsZiel.Range("C51").Value = "willl Power BI für PSI Penta, "

# --- Restore selections to match target state ---
$wsProzess.Range("D49").Select()
$wsZiel.Activate()
$wsZiel.Range("M42").Select()
